# Update "想去人数" (want-to-go count) values in the F column
# for the "展览" (sheet 1) and "全部类型" (sheet 4) sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 625
    $ws.Range("F3").Value = 3809
    $ws.Range("F4").Value = 105
    $ws.Range("F5").Value = 723
}
